# Update the "想去人数" (people interested) counts on the two sheets that
# list exhibition events: "展览" and "全部类型".
# Row 2 (丽水·新年动漫狂欢盛典): F2  328 -> 329
# Row 5 (丽水·LPJ 现实X次元动漫展): F5 284 -> 285

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 329
    $ws.Range("F5").Value = 285
}
